$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Test "
$ws.Range("C2").Value = "09:26"
$ws.Range("D2").Value = "TEst 1"
$ws.Range("F2").Value = 0

$ws.Range("B3").Value = "Test 3"
$ws.Range("C3").Value = "09:26"
$ws.Range("D3").Value = "Test 4"

$ws.Range("B4").Value = "Test 4"
$ws.Range("C4").Value = "09:26"
$ws.Range("D4").Value = "TEst 1"

$ws.Range("B5").Value = "TEst 1"
$ws.Range("C5").Value = "09:26"
$ws.Range("D5").Value = "WINNER"
$ws.Range("F5").Value = ""
